$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Save the existing Consumer_No_List value (B2) before overwriting it
$oldValue = $ws.Range("B2").Text

# Update B2 with the new consumer numbers list
$ws.Range("B2").Value = "0400005777052,0400033016199"

# Append the old value to a new row further down (row 7), matching B2's style
$ws.Range("B7").Value = $oldValue
$ws.Range("B7").NumberFormat = $ws.Range("B2").NumberFormat

# Update the active selection to reflect the new cursor position
$ws.Range("C11").Select()
